# Apply updated calibration values across the workbook's sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: start_price ---
$wsStart = $wb.Worksheets.Item("start_price")
$wsStart.Range("A2").Value = 472.9

# --- Sheet: Linear ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.05466670095456525
$wsLinear.Range("B3").Value = -0.0161811157403841
$wsLinear.Range("B4").Value = 49.6038202290494
$wsLinear.Range("B5").Value = "[1.0, 0.20984586221624715, 0.06548604999759751, 0.08375437303095808, 0.07498816548374153, 0.06295616565398102, 0.24469576128294146, 0.3451723365128955, 0.20819252268387742, 0.0634223269343312, 0.018172975394580466, 0.062779419257647, 0.06817322485850079, 0.20074977032713964, 0.3328884284719947, 0.2137265184375695, 0.05319782407448279, 0.045007146891732346, 0.048115036873010886, 0.04023490833098028]"

# --- Sheet: NonLinear ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 0.913825402772574
$wsNonLinear.Range("B4").Value = -0.2484522502493774
$wsNonLinear.Range("B5").Value = -0.06345877346708208
$wsNonLinear.Range("B6").Value = 51.54841082646607
$wsNonLinear.Range("B7").Value = 0.4831494763226852
$wsNonLinear.Range("B8").Value = -0.1351458380333136
$wsNonLinear.Range("B9").Value = 47.71880068815304
$wsNonLinear.Range("B10").Value = "[1.0, 0.20593867602897628, 0.06803458689607957, 0.08632502668054165, 0.07837768737509772, 0.06693626787079539, 0.2434636491895098, 0.3388296189808189, 0.20703699965912056, 0.0675272409679207, 0.022624236406595707, 0.06613989230059669, 0.0727270420543644, 0.19936949816551625, 0.3263130053291425, 0.21201709783687192, 0.056651691257516096, 0.04902323036608182, 0.052467453308953485, 0.04312539137351179]"
